# Horarios actualizados Linea 141 - 199
# New scrape timestamp reflected across all three sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:52:23"
$ws1.Range("A3").Value = "Total filas: 64"

# Insert the 10 newly scraped rows (ascending order of their final row number
# so each insertion lands exactly where it belongs).
$ws1.Rows.Item(33).Insert()
$ws1.Rows.Item(45).Insert()
$ws1.Rows.Item(49).Insert()
$ws1.Rows.Item(50).Insert()
$ws1.Rows.Item(56).Insert()
$ws1.Rows.Item(61).Insert()
$ws1.Rows.Item(65).Insert()
$ws1.Rows.Item(67).Insert()
$ws1.Rows.Item(68).Insert()
$ws1.Rows.Item(69).Insert()

# The re-sort by arrival time also swapped these two pre-existing rows
# (both previously arriving "07:31"), now landing at rows 52/53.
# NOTE: use .Value2 for reads here - .Value (without invocation) is not
# reliable when immediately re-assigned through a variable in this host.
$tmpA = $ws1.Cells.Item(52, 1).Value2
$tmpB = $ws1.Cells.Item(52, 2).Value2
$tmpC = $ws1.Cells.Item(52, 3).Value2
$tmpD = $ws1.Cells.Item(52, 4).Value2
$tmpE = $ws1.Cells.Item(52, 5).Value2

$ws1.Cells.Item(52, 1).Value = $ws1.Cells.Item(53, 1).Value2
$ws1.Cells.Item(52, 2).Value = $ws1.Cells.Item(53, 2).Value2
$ws1.Cells.Item(52, 3).Value = $ws1.Cells.Item(53, 3).Value2
$ws1.Cells.Item(52, 4).Value = $ws1.Cells.Item(53, 4).Value2
$ws1.Cells.Item(52, 5).Value = $ws1.Cells.Item(53, 5).Value2

$ws1.Cells.Item(53, 1).Value = $tmpA
$ws1.Cells.Item(53, 2).Value = $tmpB
$ws1.Cells.Item(53, 3).Value = $tmpC
$ws1.Cells.Item(53, 4).Value = $tmpD
$ws1.Cells.Item(53, 5).Value = $tmpE

# Fill in the freshly inserted rows with the new scrape's data.
function Set-Row141($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

Set-Row141 $ws1 33 "06:52:23" "06:52" "215C_EL PATO" 0 "LP1912"
Set-Row141 $ws1 45 "06:52:23" "07:16" "16_SANTA ANA" 24 "LP1912"
Set-Row141 $ws1 49 "06:52:23" "07:23" "10_OLMOS" 31 "LP1912"
Set-Row141 $ws1 50 "06:52:23" "07:25" "10_OLMOS" 33 "LP1912"
Set-Row141 $ws1 56 "06:52:23" "07:37" "27_EL RETIRO" 45 "LP1912"
Set-Row141 $ws1 61 "06:52:23" "08:03" "23_HERNANDEZ" 71 "LP1912"
Set-Row141 $ws1 65 "06:52:23" "08:21" "26_HERNANDEZ" 89 "LP1912"
Set-Row141 $ws1 67 "06:52:23" "08:23" "215B_EL PATO" 91 "LP1912"
Set-Row141 $ws1 68 "06:52:23" "08:27" "84_COLONIA URQUIZA-ESC 49" 95 "LP1912"
Set-Row141 $ws1 69 "06:52:23" "08:42" "81_EL PELIGRO" 110 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:52:23"
$ws2.Range("A3").Value = "Total filas: 13"

$ws2.Rows.Item(14).Insert()
$ws2.Rows.Item(18).Insert()

Set-Row141 $ws2 14 "06:52:23" "06:52" "215C_EL PATO" 0 "LP1912"
Set-Row141 $ws2 18 "06:52:23" "08:23" "215B_EL PATO" 91 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:52:23"
$ws3.Range("A3").Value = "Total filas: 13"

$ws3.Rows.Item(15).Insert()
$ws3.Rows.Item(17).Insert()
$ws3.Rows.Item(18).Insert()

Set-Row141 $ws3 15 "06:52:23" "07:38" "215A_LA PLATA" 46 "L6173"
Set-Row141 $ws3 17 "06:52:23" "08:11" "215C_LA PLATA" 79 "L6203"
Set-Row141 $ws3 18 "06:52:23" "08:40" "215A_LA PLATA" 108 "L6173"
